$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '28.031.64'
$ws.Range('E2').Value = '  +6.72%  '
Set-TextValue 'D3' '1.742.30'
$ws.Range('E3').Value = '  +5.15%  '
$ws.Range('E4').Value = '  -0.19%  '
Set-TextValue 'D5' '229.95'
$ws.Range('E5').Value = '  +4.88%  '
Set-TextValue 'D6' '0.5463'
$ws.Range('E6').Value = '  +4.17%  '
$ws.Range('E7').Value = '  -0.22%  '
Set-TextValue 'D8' '0.2791'
$ws.Range('E8').Value = '  +4.47%  '
Set-TextValue 'D9' '0.06736'
$ws.Range('E9').Value = '  +5.71%  '
Set-TextValue 'D10' '21.93'
$ws.Range('E10').Value = '  +5.85%  '
Set-TextValue 'D11' '0.07790'
$ws.Range('E11').Value = '  +0.96%  '
Set-TextValue 'D12' '4.725'
$ws.Range('E12').Value = '  +2.61%  '
Set-TextValue 'B13' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C13' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D13' '1.979.59'
$ws.Range('E13').Value = '  +4.99%  '
Set-TextValue 'B14' 'WrappedEther'
Set-TextValue 'C14' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D14' '1.717.32'
$ws.Range('E14').Value = '  +2.09%  '
Set-TextValue 'D15' '0.6029'
$ws.Range('E15').Value = '  +6.66%  '
Set-TextValue 'D16' '0.0₅8445'
$ws.Range('E16').Value = '  +2.20%  '
Set-TextValue 'D17' '69.90'
$ws.Range('E17').Value = '  +6.73%  '
Set-TextValue 'D18' '28.005.31'
$ws.Range('E18').Value = '  +6.64%  '
Set-TextValue 'D19' '228.97'
$ws.Range('E19').Value = '  +18.45%  '
Set-TextValue 'D20' '4.848'
$ws.Range('E20').Value = '  +3.23%  '
$ws.Range('E21').Value = '  -0.20%  '
Set-TextValue 'D22' '10.99'
$ws.Range('E22').Value = '  +5.20%  '
Set-TextValue 'D23' '6.286'
$ws.Range('E23').Value = '  +4.62%  '
Set-TextValue 'D24' '1.002'
Set-TextValue 'D25' '147.18'
$ws.Range('E25').Value = '  +2.69%  '
Set-TextValue 'D26' '0.1253'
$ws.Range('E26').Value = '  +4.11%  '
Set-TextValue 'D27' '7.480'
$ws.Range('E27').Value = '  +2.50%  '
Set-TextValue 'D28' '17.22'
$ws.Range('E28').Value = '  +8.01%  '
Set-TextValue 'D29' '1.617'
$ws.Range('E29').Value = '  +7.17%  '
Set-TextValue 'D30' '0.05661'
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('E31').Value = '  +3.40%  '
$ws.Range('E32').Value = '  +5.95%  '
Set-TextValue 'D33' '3.557'
$ws.Range('E33').Value = '  +6.03%  '
Set-TextValue 'D34' '1.663'
$ws.Range('E34').Value = '  +4.87%  '
Set-TextValue 'D35' '0.9873'
$ws.Range('E35').Value = '  +4.01%  '
Set-TextValue 'D36' '2.863'
$ws.Range('E36').Value = '  +1.99%  '
Set-TextValue 'D37' '2.452'
$ws.Range('E37').Value = '  +1.54%  '
Set-TextValue 'D38' '0.5960'
$ws.Range('E38').Value = '  +3.17%  '
Set-TextValue 'D39' '0.01682'
$ws.Range('E39').Value = '  +4.92%  '
Set-TextValue 'D40' '6.004'
$ws.Range('E40').Value = '  +0.44%  '
Set-TextValue 'B41' 'Maker'
Set-TextValue 'C41' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D41' '1.050.33'
$ws.Range('E41').Value = '  +3.13%  '
Set-TextValue 'B42' 'TrustWalletToken'
Set-TextValue 'C42' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D42' '0.8484'
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('E43').Value = '  -0.17%  '
Set-TextValue 'D44' '102.49'
$ws.Range('E44').Value = '  +0.58%  '
Set-TextValue 'D45' '1.884.68'
$ws.Range('E45').Value = '  +4.91%  '
$ws.Range('E46').Value = '  +11.14%  '
Set-TextValue 'D47' '60.41'
$ws.Range('E47').Value = '  +3.47%  '
Set-TextValue 'D48' '8.337'
$ws.Range('E48').Value = '  +3.55%  '
Set-TextValue 'B49' 'Frax'
Set-TextValue 'C49' 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D49' '1.012'
$ws.Range('E49').Value = '  +0.89%  '
Set-TextValue 'B50' 'Mantle'
Set-TextValue 'C50' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D50' '0.4421'
$ws.Range('E50').Value = '  +1.64%  '
Set-TextValue 'D51' '0.05314'
$ws.Range('E51').Value = '  -0.23%  '
